$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("G1").Value = "extendibile"
$ws.Range("H1").Value = "other features?"

# --- Row 2 (Bootbot) ---
$ws.Range("D2").Value = 2
$ws.Range("G2").Value = "yes"

# --- Row 3 (Dialogflow) ---
$ws.Range("D3").Value = "All"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "$0.007 per message, + other google cloud platform costs"
$ws.Range("G3").Value = "yes"
$ws.Range("H3").Value = "Designed for speech to text requests"

# --- Row 4 (Pandorabot) ---
$ws.Range("B4").Value = "Instantly"
$ws.Range("D4").Value = "All - with paid version"
$ws.Range("E4").Value = "AIML"
$ws.Range("F4").Value = "$19/month + $9/channel"
$ws.Range("G4").Value = "With paid version"
$ws.Range("H4").Value = "Industry standard"

# --- Column widths ---
# (target stored widths are 50.28515625 / 16.5703125; the engine quantizes
# ColumnWidth to steps of 1/6 + 5/6, so we feed it the values whose rounded
# result lands closest to the target)
$ws.Columns.Item(6).ColumnWidth = 49.5
$ws.Columns.Item(7).ColumnWidth = 15.6666666666667

# --- Selection ---
$ws.Range("F14").Select()
